$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.747.46'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.905.38'
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.97'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5220'
$ws.Range("E7").Value = '  +5.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3783'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07249'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.25'
$ws.Range("E10").Value = '  +3.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9039'
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07653'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.912.33'
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.451'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.14'
$ws.Range("E15").Value = '  +1.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9999'
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008709'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9996'
$ws.Range("E18").Value = '  -0.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '27.788.11'
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.53'
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.142'
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.152.26'
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.629'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.56'
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.869'
$ws.Range("E26").Value = '  +1.41%  '
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.158'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.64'
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.855'
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09045'
$ws.Range("E31").Value = '  +1.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.185'
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("E33").Value = '  +4.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.231'
$ws.Range("E34").Value = '  +0.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7815'
$ws.Range("E35").Value = '  +2.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02096'
$ws.Range("E36").Value = '  +2.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.602'
$ws.Range("E37").Value = '  +1.71%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.075'
$ws.Range("E38").Value = '  +2.97%  '
$ws.Range("E39").Value = '  -0.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5560'
$ws.Range("E40").Value = '  +1.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05289'
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.721'
$ws.Range("E42").Value = '  -2.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '114.96'
$ws.Range("E43").Value = '  +2.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.521'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1519'
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4823'
$ws.Range("E46").Value = '  +0.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.48'
$ws.Range("E47").Value = '  -1.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9993'
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.617'
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '66.81'
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06005'
$ws.Range("E51").Value = '  -0.82%  '
